# =====================================================================
# Circle Language Spec: "1. Relations Between Commands & Objects"
#
# 1. Remove the (hidden) "_GoBack" bookmark that currently sits right
#    after the second "Object relations" (Heading 4) paragraph.
# 2. Re-order / re-style the "From the original Symbol documentation"
#    (Heading 2) + "Parameter and Argument" (Heading 3) paragraph pair:
#    "Parameter and Argument" becomes the Heading 3 title, and the old
#    title text becomes an italicised, parenthesised caption right
#    below it, followed by a new blank paragraph.
# 3. Turn the last (empty) paragraph of the document into a new
#    "Misc Ideas" (Heading 3) section with a follow-up paragraph of
#    brainstormed text, and re-create the "_GoBack" bookmark inside
#    that new paragraph (mirroring where Word would leave it after the
#    last edit).
# =====================================================================

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 3 - new "Misc Ideas" section at the end of the document.
# Handled first (bottom-up) so paragraph indices used below stay valid.
# ---------------------------------------------------------------------

$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)

# The existing trailing empty paragraph becomes the new heading.
$lastPara.Range.Text = "Misc Ideas"
$lastPara.Style = "Heading 3"

# Add a fresh paragraph after the heading to hold the brainstormed text.
$headingPara = $d.Paragraphs($lastIndex)
$headingPara.Range.InsertParagraphAfter()

$ideaPara = $d.Paragraphs($lastIndex + 1)
$ideaPara.Style = "Normal"
$ideaPara.Range.Text = 'Nice phrase? "how to link objects to a command."'

# Re-create the "_GoBack" bookmark right after "Nice phrase? " (this
# also naturally splits the run there, same as the original document).
$ideaPara2 = $d.Paragraphs($lastIndex + 1)
$bmPos = $ideaPara2.Range.Start + "Nice phrase? ".Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------
# Part 2 - "From the original Symbol documentation" / "Parameter and
# Argument" reshuffle.
# ---------------------------------------------------------------------

$titlePara = $d.Paragraphs(131)
$titlePara.Range.Text = "Parameter and Argument"
$titlePara.Style = "Heading 3"

$captionPara = $d.Paragraphs(132)
$captionFull = $d.Range($captionPara.Range.Start, $captionPara.Range.End)
$captionFull.Delete()
$captionFull.InsertBefore("(From the original Symbol documentation)`r")

$captionPara2 = $d.Paragraphs(132)
$captionPara2.Range.Font.Italic = $true
$captionPara2.Range.Font.ItalicBi = $true

# New blank paragraph following the caption.
$captionPara3 = $d.Paragraphs(132)
$captionPara3.Range.InsertParagraphAfter()
$blankPara = $d.Paragraphs(133)
$blankPara.Range.Delete()

# ---------------------------------------------------------------------
# Part 1 - drop the stray "_GoBack" bookmark after the second
# "Object relations" (Heading 4) paragraph.
# ---------------------------------------------------------------------

$idx = 0
$objRelPara = $null
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -eq "Object relations`r" -and $p.Style.NameLocal -eq "Heading 4") {
        $objRelPara = $p
    }
}

$r = $objRelPara.Range
$lastCharRange = $d.Range($r.End - 2, $r.End - 1)
$savedText = $lastCharRange.Text
$lastCharRange.Delete()
$lastCharRange.InsertBefore($savedText)
